$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "id", "bug type", "species", "bug_category" columns (AC:AF),
# which shifts "predicted_bug_category" (was AG) left into column AC.
$ws.Range("AC:AF").Delete()
